$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet ("Foglio8")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "5.6 Annual summaries"

# Header row
$ws.Range("A1").Value = "row"
$ws.Range("B1").Value = "group"
$ws.Range("C1").Value = "item"
$ws.Range("D1").Value = "label"
$ws.Range("E1").Value = "unit"

# Column widths similar to other summary sheets
$ws.Columns.Item(2).ColumnWidth = 22.21875
$ws.Columns.Item(3).ColumnWidth = 22

# Data rows: group (B) / item (C) pairs
$groups = @(
    "Major power producers",
    "Major power producers",
    "Major power producers",
    "Major power producers",
    "Major power producers",
    "Major power producers",
    "Other generators",
    "Other generators",
    "Other generators",
    "Other generators",
    "All generating companies",
    "All generating companies",
    "All generating companies",
    "All generating companies",
    "All generating companies",
    "All generating companies"
)

$items = @(
    "Fuel used ",
    "Generation",
    "Used on works",
    "Supplied (gross)",
    "Used in pumping",
    "Supplied (net)",
    "Fuel used ",
    "Generation ",
    "Used on works",
    "Supplied ",
    "Fuel used",
    "Generation ",
    "Used on works",
    "Supplied (gross)",
    "Used in pumping",
    "Supplied (net)"
)

# Row index counter (A column): starts at 0
$ws.Range("A2").Value = 0
$ws.Range("A3").Formula = "=1+A2"
$ws.Range("A4:A17").Formula = "=1+A3"

for ($i = 0; $i -lt $groups.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $groups[$i]
    $ws.Cells.Item($row, 3).Value = $items[$i]
    $ws.Cells.Item($row, 5).Value = "GWh"
}

$ws.Range("D2").Formula = '=_xlfn.CONCAT(B2," ",C2)'
$ws.Range("D3:D17").Formula = '=_xlfn.CONCAT(B3," ",C3)'

# Sheet view tweaks to match the target: zoomed, scrolled, tab selected
$ws.Application.ActiveWindow.Zoom = 190
$ws.Range("A18").Select()
